$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving exact string (avoids Excel
# auto-converting numeric-looking strings to numbers and avoids leaving a
# lingering custom cell style behind).
function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.229.96'
$ws.Range("E2").Value = '  -0.14%  '

Set-TextValue $ws.Range("D3") '1.856.42'
$ws.Range("E3").Value = '  -0.28%  '

Set-TextValue $ws.Range("D4") '0.9997'
$ws.Range("E4").Value = '  +0.01%  '

Set-TextValue $ws.Range("D5") '241.03'
$ws.Range("E5").Value = '  -0.61%  '

Set-TextValue $ws.Range("D6") '0.6983'
$ws.Range("E6").Value = '  -0.91%  '

Set-TextValue $ws.Range("D7") '0.9998'
$ws.Range("E7").Value = '  +0.01%  '

Set-TextValue $ws.Range("D8") '0.07782'
$ws.Range("E8").Value = '  -0.24%  '

Set-TextValue $ws.Range("D9") '0.3075'
$ws.Range("E9").Value = '  -2.53%  '

Set-TextValue $ws.Range("D10") '23.74'
$ws.Range("E10").Value = '  -2.17%  '

Set-TextValue $ws.Range("D11") '0.07793'
$ws.Range("E11").Value = '  -2.64%  '

Set-TextValue $ws.Range("D12") '1.866.13'
$ws.Range("E12").Value = '  +0.26%  '

Set-TextValue $ws.Range("D13") '5.110'
$ws.Range("E13").Value = '  -1.53%  '

Set-TextValue $ws.Range("D14") '92.28'
$ws.Range("E14").Value = '  -1.96%  '

Set-TextValue $ws.Range("D15") '0.6878'
$ws.Range("E15").Value = '  -1.38%  '

Set-TextValue $ws.Range("D16") '6.529'
$ws.Range("E16").Value = '  +1.53%  '

Set-TextValue $ws.Range("D17") '0.000008463'
$ws.Range("E17").Value = '  +2.14%  '

Set-TextValue $ws.Range("D18") '29.222.54'
$ws.Range("E18").Value = '  -0.18%  '

Set-TextValue $ws.Range("D19") '248.10'
$ws.Range("E19").Value = '  -1.95%  '

Set-TextValue $ws.Range("D20") '2.108.27'
$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("E21").Value = '  -2.30%  '

Set-TextValue $ws.Range("D22") '1.0000'
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("E23").Value = '  -0.15%  '

Set-TextValue $ws.Range("D24") '0.9998'
$ws.Range("E24").Value = '  +0.00%  '

Set-TextValue $ws.Range("D25") '0.1504'
$ws.Range("E25").Value = '  -4.04%  '

Set-TextValue $ws.Range("D26") '161.65'
$ws.Range("E26").Value = '  +0.91%  '

Set-TextValue $ws.Range("D27") '8.862'
$ws.Range("E27").Value = '  -1.52%  '

Set-TextValue $ws.Range("D28") '18.51'
$ws.Range("E28").Value = '  -1.94%  '

$ws.Range("E29").Value = '  +3.76%  '

Set-TextValue $ws.Range("D30") '4.260'
$ws.Range("E30").Value = '  -1.26%  '

Set-TextValue $ws.Range("D31") '4.210'
$ws.Range("E31").Value = '  -1.42%  '

Set-TextValue $ws.Range("D32") '1.196'
$ws.Range("E32").Value = '  -0.93%  '

$ws.Range("E33").Value = '  -1.15%  '

Set-TextValue $ws.Range("D34") '0.7650'
$ws.Range("E34").Value = '  +1.90%  '

Set-TextValue $ws.Range("D35") '1.848'
$ws.Range("E35").Value = '  -2.17%  '

Set-TextValue $ws.Range("D36") '1.170'
$ws.Range("E36").Value = '  +0.93%  '

Set-TextValue $ws.Range("D37") '2.709'
$ws.Range("E37").Value = '  -0.14%  '

Set-TextValue $ws.Range("D38") '0.01864'
$ws.Range("E38").Value = '  -0.30%  '

Set-TextValue $ws.Range("D39") '1.226.14'
$ws.Range("E39").Value = '  -1.93%  '

$ws.Range("E40").Value = '  -0.36%  '

Set-TextValue $ws.Range("D41") '0.9006'
$ws.Range("E41").Value = '  +0.15%  '

Set-TextValue $ws.Range("D42") '109.33'
$ws.Range("E42").Value = '  -1.56%  '

$ws.Range("E43").Value = '  +0.00%  '

Set-TextValue $ws.Range("D44") '5.567'
$ws.Range("E44").Value = '  -9.27%  '

Set-TextValue $ws.Range("D45") '2.006.11'
$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("E46").Value = '  -3.35%  '

Set-TextValue $ws.Range("D47") '65.40'
$ws.Range("E47").Value = '  -7.49%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D48") '0.5184'
$ws.Range("E48").Value = '  -0.14%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D49") '9.556'
$ws.Range("E49").Value = '  +0.76%  '

$ws.Range("E50").Value = '  -1.99%  '

Set-TextValue $ws.Range("D51") '7.044'
$ws.Range("E51").Value = '  +0.18%  '
